# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.581.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.08%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.410.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.91%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.10%  "

# Row 9
$ws.Range("E9").Value = "  +4.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.361"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.45%  "

# Row 12
$ws.Range("E12").Value = "  -2.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.844.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.495.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.07%  "

# Row 16
$ws.Range("E16").Value = "  +4.29%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.396.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.55%  "

# Row 19
$ws.Range("E19").Value = "  +4.46%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.65%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.60%  "

# Row 22
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.75%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.169"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  -1.67%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0788"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.84%  "

# Row 29
$ws.Range("E29").Value = "  +2.44%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.55%  "

# Row 34
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("E35").Value = "  +5.18%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.53%  "

# Row 37
$ws.Range("E37").Value = "  +0.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.420"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "306.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.11%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "142.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0964"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.04%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.571"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.400"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "

# Row 51
$ws.Range("E51").Value = "  +4.76%  "
